{"js": "// Update the \"Gi\u1ea5y \u0111\u1ec1 ngh\u1ecb thanh to\u00e1n\" template: replace the requester's\n// name and clear out the autofilled values for \u0110\u01a1n v\u1ecb, N\u1ed9i dung thanh\n// to\u00e1n, M\u00e3 s\u1ed1 \u0111\u1ec1 t\u00e0i, S\u1ed1 ti\u1ec1n, and Vi\u1ebft b\u1eb1ng ch\u1eef (keeping their labels).\n\nconst replacements = [\n  {\n    find: \"H\u1ecd, t\u00ean ng\u01b0\u1eddi \u0111\u1ec1 ngh\u1ecb thanh to\u00e1n: L\u00ea Anh Ti\u1ebfn\",\n    replace: \"H\u1ecd, t\u00ean ng\u01b0\u1eddi \u0111\u1ec1 ngh\u1ecb thanh to\u00e1n: \u0110\u1ed7 Thanh B\u00ecnh\",\n  },\n  {\n    find: \"\u0110\u01a1n v\u1ecb: Khoa C\u00f4ng Ngh\u1ec7 Th\u00f4ng Tin\",\n    replace: \"\u0110\u01a1n v\u1ecb: \",\n  },\n  {\n    find:\n      \"N\u1ed9i dung thanh to\u00e1n: Thanh to\u00e1n kinh ph\u00ed c\u1ee7a \u0111\u1ec3 t\u00e0i c\u1ea5p c\u01a1 s\u1edf Nghi\u00ean c\u1ee9u ki\u1ebfn tr\u00fac chip x\u1eed l\u00fd m\u1eadt m\u00e3 theo ti\u00eau chu\u1ea9n Trusted Platform Module 2.0 (TPM 2.0) c\u1ee7a Trusted Computing Group (TCG)\\\" M\u00e3 s\u1ed1 \u0111\u1ec3 t\u00e0i 19/2023/CS do L\u00ea Anh Ti\u1ebfn l\u00e0m ch\u1ee7 nhi\u1ec7m\",\n    replace: \"N\u1ed9i dung thanh to\u00e1n: \",\n  },\n  {\n    find: \"M\u00e3 s\u1ed1 \u0111\u1ec1 t\u00e0i: 19/2023/CS do L\u00ea Anh Ti\u1ebfn l\u00e0m ch\u1ee7 nhi\u1ec7m\",\n    replace: \"M\u00e3 s\u1ed1 \u0111\u1ec1 t\u00e0i: \",\n  },\n  {\n    find: \"S\u1ed1 ti\u1ec1n: 50.000.000 \u0111\u1ed3ng\",\n    replace: \"S\u1ed1 ti\u1ec1n: \",\n  },\n  {\n    find: \"Vi\u1ebft b\u1eb1ng ch\u1eef: N\u0103m m\u01b0\u01a1i tri\u1ec7u \u0111\u1ed3ng ch\u1eb3n\",\n    replace: \"Vi\u1ebft b\u1eb1ng ch\u1eef: \",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    continue;\n  }\n\n  results.items[0].insertText(replace, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the \"Gi\u1ea5y \u0111\u1ec1 ngh\u1ecb thanh to\u00e1n\" template: replace the requester's\n# name and clear out the autofilled values for \u0110\u01a1n v\u1ecb, N\u1ed9i dung thanh\n# to\u00e1n, M\u00e3 s\u1ed1 \u0111\u1ec1 t\u00e0i, S\u1ed1 ti\u1ec1n, and Vi\u1ebft b\u1eb1ng ch\u1eef (keeping their labels).\n\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n\n$replacements = @(\n    @{\n        Find    = \"H\u1ecd, t\u00ean ng\u01b0\u1eddi \u0111\u1ec1 ngh\u1ecb thanh to\u00e1n: L\u00ea Anh Ti\u1ebfn\"\n        Replace = \"H\u1ecd, t\u00ean ng\u01b0\u1eddi \u0111\u1ec1 ngh\u1ecb thanh to\u00e1n: \u0110\u1ed7 Thanh B\u00ecnh\"\n    },\n    @{\n        Find    = \"\u0110\u01a1n v\u1ecb: Khoa C\u00f4ng Ngh\u1ec7 Th\u00f4ng Tin\"\n        Replace = \"\u0110\u01a1n v\u1ecb: \"\n    },\n    @{\n        Find    = \"N\u1ed9i dung thanh to\u00e1n: Thanh to\u00e1n kinh ph\u00ed c\u1ee7a \u0111\u1ec3 t\u00e0i c\u1ea5p c\u01a1 s\u1edf Nghi\u00ean c\u1ee9u ki\u1ebfn tr\u00fac chip x\u1eed l\u00fd m\u1eadt m\u00e3 theo ti\u00eau chu\u1ea9n Trusted Platform Module 2.0 (TPM 2.0) c\u1ee7a Trusted Computing Group (TCG)`\" M\u00e3 s\u1ed1 \u0111\u1ec3 t\u00e0i 19/2023/CS do L\u00ea Anh Ti\u1ebfn l\u00e0m ch\u1ee7 nhi\u1ec7m\"\n        Replace = \"N\u1ed9i dung thanh to\u00e1n: \"\n    },\n    @{\n        Find    = \"M\u00e3 s\u1ed1 \u0111\u1ec1 t\u00e0i: 19/2023/CS do L\u00ea Anh Ti\u1ebfn l\u00e0m ch\u1ee7 nhi\u1ec7m\"\n        Replace = \"M\u00e3 s\u1ed1 \u0111\u1ec1 t\u00e0i: \"\n    },\n    @{\n        Find    = \"S\u1ed1 ti\u1ec1n: 50.000.000 \u0111\u1ed3ng\"\n        Replace = \"S\u1ed1 ti\u1ec1n: \"\n    },\n    @{\n        Find    = \"Vi\u1ebft b\u1eb1ng ch\u1eef: N\u0103m m\u01b0\u01a1i tri\u1ec7u \u0111\u1ed3ng ch\u1eb3n\"\n        Replace = \"Vi\u1ebft b\u1eb1ng ch\u1eef: \"\n    }\n)\n\nforeach ($item in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $item.Find\n    $find.Replacement.Text = $item.Replace\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, $null, $false, $null, $wdReplaceAll)\n}\n"}
